$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Round the coordinate values in row 3 to whole numbers
$ws.Range("Q3").Value = 501495
$ws.Range("R3").Value = 7036929

# Clear the time cells (Starttid / Sluttid) for row 3, leaving Slutdatum (AA3) untouched
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
